# Update "想去人数" (interest count) figures in the 展览 and 全部类型 sheets
# to the latest values generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3879
$ws1.Range("F4").Value = 2304
$ws1.Range("F5").Value = 456
$ws1.Range("F7").Value = 24
$ws1.Range("F10").Value = 108
$ws1.Range("F11").Value = 1434
$ws1.Range("F12").Value = 254
$ws1.Range("F13").Value = 2543
$ws1.Range("F14").Value = 179

# Sheet "全部类型" (all categories)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3879
$ws4.Range("F4").Value = 2304
$ws4.Range("F5").Value = 456
$ws4.Range("F7").Value = 24
$ws4.Range("F11").Value = 108
$ws4.Range("F14").Value = 1434
$ws4.Range("F15").Value = 254
$ws4.Range("F16").Value = 2543
$ws4.Range("F17").Value = 179
